{"js": "// Replace each old text with its new text using Body.search + insertText(replace).\n// Using matchCase + exact whole strings so each search pinpoints only the\n// intended run(s); every old string below is unique in the document except\n// the title, which intentionally appears twice (heading + closing bold\n// recap) and gets the same replacement both times.\nconst replacements = [\n  [\n    \"Play Dead or Alive Slot Game for Free\",\n    \"Play Dead or Alive - Free Wild West Slot Game\",\n  ],\n  [\n    \"Sticky Win feature adds extra excitement\",\n    \"Immersive gameplay with Wild West theme\",\n  ],\n  [\n    \"Top-notch graphics and sound effects\",\n    \"Sticky Win feature adds excitement\",\n  ],\n  [\n    \"Immersive Western-inspired atmosphere\",\n    \"High-quality graphics and sound effects\",\n  ],\n  [\n    \"Range of betting options to suit all budgets\",\n    \"Chance to win big with free spins and Sticky Win feature\",\n  ],\n  [\n    \"Limited variety of bonus features\",\n    \"Limited betting options\",\n  ],\n  [\n    \"High variance may not appeal to all players\",\n    \"Autoplay feature may not appeal to all players\",\n  ],\n  [\n    \"Read our review of Dead or Alive slot game by NetEnt and play for free. Enjoy the immersive Western-themed atmosphere and Sticky Win feature for big wins.\",\n    \"Read our review of Dead or Alive, the free Wild West slot game with immersive gameplay.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit: rewrite the title, recap\n# heading/summary, and the pros/cons bullet lists via Find & Replace.\n# wdReplaceAll (2) so the title - which legitimately appears twice\n# (H1 + closing bold recap) - gets updated in both spots.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Dead or Alive Slot Game for Free\", \"Play Dead or Alive - Free Wild West Slot Game\"),\n    @(\"Sticky Win feature adds extra excitement\", \"Immersive gameplay with Wild West theme\"),\n    @(\"Top-notch graphics and sound effects\", \"Sticky Win feature adds excitement\"),\n    @(\"Immersive Western-inspired atmosphere\", \"High-quality graphics and sound effects\"),\n    @(\"Range of betting options to suit all budgets\", \"Chance to win big with free spins and Sticky Win feature\"),\n    @(\"Limited variety of bonus features\", \"Limited betting options\"),\n    @(\"High variance may not appeal to all players\", \"Autoplay feature may not appeal to all players\"),\n    @(\"Read our review of Dead or Alive slot game by NetEnt and play for free. Enjoy the immersive Western-themed atmosphere and Sticky Win feature for big wins.\", \"Read our review of Dead or Alive, the free Wild West slot game with immersive gameplay.\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith, Replace(wdReplaceAll=2)\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
